# Atualização 06/07 - construtor do TimeSheetManager: verifica/cria a planilha
# "PontoEletrônico.xlsx" e grava o cabeçalho da tabela de ponto eletrônico.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renomeia a aba (antes "Time sheet register") para o nome do arquivo.
$ws.Name = "PontoEletrônico.xlsx"

# Escreve o cabeçalho da planilha de ponto eletrônico na linha 1.
$ws.Range("A1").Value = "DATA"
$ws.Range("B1").Value = "ENTRADA"
$ws.Range("C1").Value = "INTERVALO"
$ws.Range("D1").Value = "RETORNO INTERVALO"
$ws.Range("E1").Value = "SAÍDA"

# Ajusta a largura das colunas ao conteúdo (best fit), como o Excel faz
# automaticamente ao criar a planilha pela primeira vez.
$ws.Range("A1:E1").EntireColumn.AutoFit()
